$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Rename header labels: "_old" -> "_FV2210", "_new" -> "_FV2304" ---
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J carry the "_old" -> "_FV2210" headers (col 1-10)
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = "$($baseNames[$i])_FV2210"
}

# Column K (11) is "diff" - unchanged

# Columns L-U carry the "_new" -> "_FV2304" headers (col 12-21)
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = "$($baseNames[$i])_FV2304"
}

# --- 2. Turn the used range into an Excel Table (ListObject) ---
$tableRange = $ws.Range("A1:U83")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# --- 3. Freeze the header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
